$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "wav.اتصل 26"
$ws.Range("B2").Value = "Call_1.wav"
$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "0:01:33"
$ws.Range("E2").Value = "'5861"
$ws.Range("F2").Value = "سزدفغ"
$ws.Range("G2").Value = "'12-12-1212"
$ws.Range("H2").Value = "الإنكليزية"
$ws.Range("I2").Value = "سالب"
$ws.Range("J2").Value = "'0.0"
$ws.Range("K2").Value = "'0.0"

# Row 3
$ws.Range("A3").Value = "wav.اتصل 27"
$ws.Range("B3").Value = "Call_2.wav"
$ws.Range("C3").Value = "'2"
$ws.Range("D3").Value = "0:01:38"
$ws.Range("E3").Value = "'6147"
$ws.Range("F3").Value = "سزدفغ"
$ws.Range("G3").Value = "'12-12-1212"
$ws.Range("H3").Value = "الإنكليزية"
$ws.Range("I3").Value = "موجب"
$ws.Range("J3").Value = "'0.8571428571428571"
$ws.Range("K3").Value = "'0.5"

# Row 4
$ws.Range("A4").Value = "wav.اتصل 28"
$ws.Range("B4").Value = "Call_3.wav"
$ws.Range("C4").Value = "'3"
$ws.Range("D4").Value = "0:01:02"
$ws.Range("E4").Value = "'3885"
$ws.Range("F4").Value = "سزدفغ"
$ws.Range("G4").Value = "'12-12-1212"

# Row 5
$ws.Range("A5").Value = "wav.اتصل 29"
$ws.Range("B5").Value = "Call_4.wav"
$ws.Range("C5").Value = "'4"
$ws.Range("D5").Value = "0:01:29"
$ws.Range("E5").Value = "'5622"
$ws.Range("F5").Value = "سزدفغ"
$ws.Range("G5").Value = "'12-12-1212"
$ws.Range("H5").Value = "الإنكليزية"
$ws.Range("I5").Value = "سالب"
$ws.Range("J5").Value = "'0.0"
$ws.Range("K5").Value = "'1.0"

# Row 6
$ws.Range("A6").Value = "wav.اتصل 30"
$ws.Range("B6").Value = "Call_5.wav"
$ws.Range("C6").Value = "'5"
$ws.Range("D6").Value = "0:01:30"
$ws.Range("E6").Value = "'5640"
$ws.Range("F6").Value = "سزدفغ"
$ws.Range("G6").Value = "'12-12-1212"
$ws.Range("H6").Value = "الإنكليزية"
$ws.Range("I6").Value = "سالب"
$ws.Range("J6").Value = "'0.0"
$ws.Range("K6").Value = "'0.0"
